{"js": "// Replace each three-digit-by-one-digit multiplication equation in the\n// document's table cells with its updated equivalent, per the commit diff.\n// Each old equation string is unique in the document, so an exact,\n// case-sensitive search-and-replace for each pair is safe and unambiguous.\nconst replacements = [\n  [\"748\u00d77=5236\", \"229\u00d73=687\"],\n  [\"225\u00d78=1800\", \"401\u00d74=1604\"],\n  [\"798\u00d79=7182\", \"490\u00d78=3920\"],\n  [\"391\u00d74=1564\", \"419\u00d76=2514\"],\n  [\"854\u00d74=3416\", \"716\u00d74=2864\"],\n  [\"145\u00d73=435\", \"178\u00d79=1602\"],\n  [\"292\u00d79=2628\", \"345\u00d77=2415\"],\n  [\"219\u00d79=1971\", \"486\u00d78=3888\"],\n  [\"346\u00d76=2076\", \"226\u00d73=678\"],\n  [\"381\u00d73=1143\", \"843\u00d73=2529\"],\n  [\"113\u00d76=678\", \"971\u00d78=7768\"],\n  [\"392\u00d74=1568\", \"931\u00d79=8379\"],\n  [\"523\u00d76=3138\", \"913\u00d73=2739\"],\n  [\"690\u00d74=2760\", \"925\u00d72=1850\"],\n  [\"951\u00d77=6657\", \"284\u00d73=852\"],\n  [\"234\u00d78=1872\", \"623\u00d78=4984\"],\n  [\"236\u00d78=1888\", \"803\u00d79=7227\"],\n  [\"140\u00d77=980\", \"920\u00d77=6440\"],\n  [\"746\u00d72=1492\", \"540\u00d73=1620\"],\n  [\"894\u00d76=5364\", \"820\u00d72=1640\"],\n  [\"594\u00d75=2970\", \"386\u00d74=1544\"],\n  [\"399\u00d78=3192\", \"878\u00d75=4390\"],\n  [\"661\u00d73=1983\", \"514\u00d78=4112\"],\n  [\"168\u00d73=504\", \"424\u00d77=2968\"],\n  [\"364\u00d72=728\", \"216\u00d79=1944\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "# Replace each three-digit-by-one-digit multiplication equation in the\n# document's table cells with its updated equivalent, per the commit diff.\n# Each old equation string is unique in the document, so an exact,\n# case-sensitive Find/Replace for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"748\u00d77=5236\"; New = \"229\u00d73=687\" },\n    @{ Old = \"225\u00d78=1800\"; New = \"401\u00d74=1604\" },\n    @{ Old = \"798\u00d79=7182\"; New = \"490\u00d78=3920\" },\n    @{ Old = \"391\u00d74=1564\"; New = \"419\u00d76=2514\" },\n    @{ Old = \"854\u00d74=3416\"; New = \"716\u00d74=2864\" },\n    @{ Old = \"145\u00d73=435\"; New = \"178\u00d79=1602\" },\n    @{ Old = \"292\u00d79=2628\"; New = \"345\u00d77=2415\" },\n    @{ Old = \"219\u00d79=1971\"; New = \"486\u00d78=3888\" },\n    @{ Old = \"346\u00d76=2076\"; New = \"226\u00d73=678\" },\n    @{ Old = \"381\u00d73=1143\"; New = \"843\u00d73=2529\" },\n    @{ Old = \"113\u00d76=678\"; New = \"971\u00d78=7768\" },\n    @{ Old = \"392\u00d74=1568\"; New = \"931\u00d79=8379\" },\n    @{ Old = \"523\u00d76=3138\"; New = \"913\u00d73=2739\" },\n    @{ Old = \"690\u00d74=2760\"; New = \"925\u00d72=1850\" },\n    @{ Old = \"951\u00d77=6657\"; New = \"284\u00d73=852\" },\n    @{ Old = \"234\u00d78=1872\"; New = \"623\u00d78=4984\" },\n    @{ Old = \"236\u00d78=1888\"; New = \"803\u00d79=7227\" },\n    @{ Old = \"140\u00d77=980\"; New = \"920\u00d77=6440\" },\n    @{ Old = \"746\u00d72=1492\"; New = \"540\u00d73=1620\" },\n    @{ Old = \"894\u00d76=5364\"; New = \"820\u00d72=1640\" },\n    @{ Old = \"594\u00d75=2970\"; New = \"386\u00d74=1544\" },\n    @{ Old = \"399\u00d78=3192\"; New = \"878\u00d75=4390\" },\n    @{ Old = \"661\u00d73=1983\"; New = \"514\u00d78=4112\" },\n    @{ Old = \"168\u00d73=504\"; New = \"424\u00d77=2968\" },\n    @{ Old = \"364\u00d72=728\"; New = \"216\u00d79=1944\" }\n)\n\n$wdReplaceAll = 2\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, $wdReplaceAll) | Out-Null\n}\n"}
